# Scheduled-runner refresh of leve profit calcs (currentAveragePrice /
# currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ) across the
# crafting-job sheets, pulling updated market-board prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1213.1818
$ws.Range("I115").Value = 1213.1818
$ws.Range("K115").Value = 3639.5454
$ws.Range("M115").Value = -2072.5454

$ws.Range("H134").Value = 54870
$ws.Range("J134").Value = 54870
$ws.Range("L134").Value = 54870
$ws.Range("N134").Value = -65010

$ws.Range("H137").Value = 1906996.6
$ws.Range("I137").Value = 2646889.8
$ws.Range("K137").Value = 7940669.399999999
$ws.Range("M137").Value = -7938119.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4987.566
$ws.Range("I32").Value = 4021.851
$ws.Range("J32").Value = 12552.333
$ws.Range("K32").Value = 4021.851
$ws.Range("L32").Value = 12552.333
$ws.Range("M32").Value = -3734.851
$ws.Range("N32").Value = -13126.333

$ws.Range("H61").Value = 1991.8334
$ws.Range("I61").Value = 1386.8572
$ws.Range("K61").Value = 1386.8572
$ws.Range("M61").Value = -1174.8572

$ws.Range("H74").Value = 11074.667
$ws.Range("I74").Value = 15871.6
$ws.Range("J74").Value = 5078.5
$ws.Range("K74").Value = 15871.6
$ws.Range("L74").Value = 5078.5
$ws.Range("M74").Value = -14997.6
$ws.Range("N74").Value = -6826.5

$ws.Range("H77").Value = 11074.667
$ws.Range("I77").Value = 15871.6
$ws.Range("J77").Value = 5078.5
$ws.Range("K77").Value = 79358
$ws.Range("L77").Value = 25392.5
$ws.Range("M77").Value = -74990
$ws.Range("N77").Value = -34128.5

$ws.Range("H122").Value = 4278
$ws.Range("I122").Value = 1278.8572
$ws.Range("J122").Value = 14775
$ws.Range("K122").Value = 3836.5716
$ws.Range("L122").Value = 44325
$ws.Range("M122").Value = -1386.5716
$ws.Range("N122").Value = -49225

$ws.Range("H132").Value = 1682.2222
$ws.Range("I132").Value = 837.4074000000001
$ws.Range("K132").Value = 2512.2222
$ws.Range("M132").Value = 17.77779999999984

$ws.Range("H136").Value = 1991.8334
$ws.Range("I136").Value = 1386.8572
$ws.Range("K136").Value = 4160.571599999999
$ws.Range("M136").Value = -1610.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 468.35
$ws.Range("I64").Value = 505.75
$ws.Range("J64").Value = 443.41666
$ws.Range("K64").Value = 505.75
$ws.Range("L64").Value = 443.41666
$ws.Range("M64").Value = -280.75
$ws.Range("N64").Value = -893.41666

$ws.Range("H67").Value = 468.35
$ws.Range("I67").Value = 505.75
$ws.Range("J67").Value = 443.41666
$ws.Range("K67").Value = 505.75
$ws.Range("L67").Value = 443.41666
$ws.Range("M67").Value = 274.25
$ws.Range("N67").Value = -2003.41666

$ws.Range("H134").Value = 1975.079
$ws.Range("I134").Value = 1635.0938
$ws.Range("J134").Value = 3788.3333
$ws.Range("K134").Value = 4905.2814
$ws.Range("L134").Value = 11364.9999
$ws.Range("M134").Value = -2370.2814
$ws.Range("N134").Value = -16434.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3349.889
$ws.Range("I31").Value = 1152.9166
$ws.Range("J31").Value = 5107.467
$ws.Range("K31").Value = 1152.9166
$ws.Range("L31").Value = 5107.467
$ws.Range("M31").Value = -857.9166
$ws.Range("N31").Value = -5697.467

$ws.Range("H34").Value = 3349.889
$ws.Range("I34").Value = 1152.9166
$ws.Range("J34").Value = 5107.467
$ws.Range("K34").Value = 1152.9166
$ws.Range("L34").Value = 5107.467
$ws.Range("M34").Value = -950.9166
$ws.Range("N34").Value = -5511.467

$ws.Range("H58").Value = 3396
$ws.Range("J58").Value = 8750
$ws.Range("L58").Value = 8750
$ws.Range("N58").Value = -9156

$ws.Range("H132").Value = 2166.8462
$ws.Range("I132").Value = 1386.8214
$ws.Range("K132").Value = 4160.4642
$ws.Range("M132").Value = -1630.4642

$ws.Range("H134").Value = 1523.5
$ws.Range("I134").Value = 844.04
$ws.Range("J134").Value = 3410.889
$ws.Range("K134").Value = 2532.12
$ws.Range("L134").Value = 10232.667
$ws.Range("M134").Value = 2.880000000000109
$ws.Range("N134").Value = -15302.667

$ws.Range("H136").Value = 3396
$ws.Range("J136").Value = 8750
$ws.Range("L136").Value = 26250
$ws.Range("N136").Value = -31350

$ws.Range("H137").Value = 45514.285
$ws.Range("J137").Value = 45514.285
$ws.Range("L137").Value = 45514.285
$ws.Range("N137").Value = -55714.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 249.14285
$ws.Range("I14").Value = 249.14285
$ws.Range("K14").Value = 747.4285500000001
$ws.Range("M14").Value = -574.4285500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1865.5714
$ws.Range("I102").Value = 1219.8276
$ws.Range("K102").Value = 1219.8276
$ws.Range("M102").Value = 402.1723999999999

$ws.Range("H107").Value = 7937003
$ws.Range("I107").Value = 293.2
$ws.Range("J107").Value = 27778778
$ws.Range("K107").Value = 293.2
$ws.Range("L107").Value = 27778778
$ws.Range("M107").Value = 1626.8
$ws.Range("N107").Value = -27782618

$ws.Range("H132").Value = 2837.0356
$ws.Range("I132").Value = 1289.6471
$ws.Range("J132").Value = 5228.4546
$ws.Range("K132").Value = 3868.9413
$ws.Range("L132").Value = 15685.3638
$ws.Range("M132").Value = -1338.9413
$ws.Range("N132").Value = -20745.3638

$ws.Range("H138").Value = 39973.332
$ws.Range("J138").Value = 39973.332
$ws.Range("L138").Value = 39973.332
$ws.Range("N138").Value = -50253.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 30000
$ws.Range("I39").Value = 30000
$ws.Range("K39").Value = 30000
$ws.Range("M39").Value = -29540

$ws.Range("H61").Value = 2613.9
$ws.Range("J61").Value = 2300.5
$ws.Range("L61").Value = 2300.5
$ws.Range("N61").Value = -2704.5

$ws.Range("H113").Value = 2613.9
$ws.Range("J113").Value = 2300.5
$ws.Range("L113").Value = 2300.5
$ws.Range("N113").Value = -6640.5

$ws.Range("H132").Value = 5174.6875
$ws.Range("I132").Value = 2799.6
$ws.Range("J132").Value = 6254.273
$ws.Range("K132").Value = 8398.799999999999
$ws.Range("L132").Value = 18762.819
$ws.Range("M132").Value = -5868.799999999999
$ws.Range("N132").Value = -23822.819

$ws.Range("H136").Value = 4049.6924
$ws.Range("I136").Value = 1213.7142
$ws.Range("K136").Value = 3641.1426
$ws.Range("M136").Value = -1091.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11496578
$ws.Range("I132").Value = 1645.45
$ws.Range("J132").Value = 37040870
$ws.Range("K132").Value = 4936.35
$ws.Range("L132").Value = 111122610
$ws.Range("M132").Value = -2406.35
$ws.Range("N132").Value = -111127670

$ws.Range("H136").Value = 4484.7
$ws.Range("I136").Value = 2736
$ws.Range("J136").Value = 5821.9414
$ws.Range("K136").Value = 8208
$ws.Range("L136").Value = 17465.8242
$ws.Range("M136").Value = -5658
$ws.Range("N136").Value = -22565.8242
